$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "AK2" = 0.06926169836713743
    "AK4" = 0.5133108883317492
    "AK5" = 0.1127443974384373
    "AK7" = 0.1589923842336786
    "AK8" = 1.11572727492984
    "AK9" = 0.05187550504318959
    "AK10" = 0.1859988056220597
    "AK11" = 0.9269315712255368
    "AK12" = 0.2191730765526122
    "AK13" = -0.4432130382515084
    "AK15" = -0.2442516445297177
    "AK19" = 0.8750846915268385
    "AK20" = -0.1702074627298488
    "AK21" = -0.003469305910314809
    "AK22" = 0.4934518945486912
    "AK24" = -1.972166678823926
    "AK25" = -1.087972011552671
    "AK26" = -0.1960190910187525
    "AK28" = -1.240597241501376
    "AK29" = -0.2991128741876304
    "AK30" = 0.4085145584359414
    "AK31" = 0.6241274801307243
    "AK32" = 0.008275432386746408
    "AK33" = 0.09354480023006469
    "AK34" = -0.8228931488316266
    "AK36" = 0.7008667685577774
    "AK37" = -0.1259468601743914
    "AK38" = -0.4235853852459345
    "AK39" = 0.5328945785188893
    "AK40" = -2.04010234783437
    "AK41" = -2.174523349698429
    "AK43" = 0.4204175717376188
    "AK46" = -0.7970124593326228
    "AK47" = -0.3513754254930812
    "AK48" = -1.688852682012541
    "AK49" = 0.3397165551053185
    "AK50" = -1.663711016068243
    "AK52" = -0.8199392948410563
    "AK53" = -0.1849290071183466
    "AK54" = -0.00330345636574594
    "AK55" = 0.004969667249321938
    "AK56" = -0.888609120227366
    "AK58" = 0.3092426912330473
    "AK59" = -3.080117332062254
    "AK60" = -3.109619734100242
    "AK64" = 0.699857151574513
    "AK65" = 0.01509391539316502
    "AK66" = 0.3279433050087912
    "AK67" = 0.9623864375566633
    "AK68" = 0.03627448048536463
    "AK71" = 0.6969491669194033
    "AK72" = -2.431102857870564
    "AK75" = 1.013308028303126
    "AK76" = -0.07519723739709405
    "AK78" = 0.4704006628639727
    "AK81" = 0.02871676232132448
    "AK82" = -0.07779048644755447
    "AK85" = -0.7647301797721917
    "AK86" = -0.2419153657985639
    "AK87" = 0.9093219155000927
    "AK88" = -0.927902055791858
    "AK89" = 0.7753211642599774
    "AK91" = 0.4472669909700842
    "AK92" = 0.1599823474581356
    "AK95" = -0.2103242957121566
    "AK96" = -0.1862750507847621
    "AK98" = -0.8313320885949088
    "AK99" = -0.07130532420167862
    "AK101" = -1.7706891639748
    "AK102" = -1.277990077086041
    "AK107" = 0.6519024537176685
    "AK110" = -0.8083417276770196
    "AK111" = -1.796775905961965
    "AK112" = 0.4808921919420336
    "AK113" = 0.281105990775068
    "AK115" = 0.1676729881882217
    "AK119" = -0.3078273142744736
    "AK120" = -0.394347883671995
    "AK121" = 0.251187010995414
    "AK122" = 0.1708853687452332
    "AK123" = -0.04708365034770719
    "AK125" = -1.728313629093525
    "AK126" = -0.06403046834052387
    "AK127" = -0.3069045578505933
    "AK129" = 0.2404941047732678
    "AK130" = 0.3147184981613569
    "AK131" = 0.6481522556954038
    "AK132" = -0.02673093086952361
    "AK133" = 0.8948848096238828
    "AK134" = 0.3928838871229469
    "AK137" = 0.4140787321128346
    "AK138" = -2.683715764195747
    "AK139" = -1.206569792703051
    "AK140" = 0.1425565368295368
    "AK141" = 0.1302388615197772
    "AK142" = -0.2745717552201338
    "AK143" = -0.7717199703195817
    "AK144" = -0.1491172619952754
    "AK145" = -0.03405750164785184
    "AK146" = 0.06436714771638037
    "AK149" = -0.2891514461210318
    "AK150" = 0.04126160151541105
    "AK151" = -0.04184343050141345
    "AK153" = -0.1654015016949656
    "AK155" = 0.1924262620971286
    "AK157" = -0.944200788308289
    "AK159" = 0.5321797191327672
    "AK161" = -0.6548047897723207
    "AK166" = 0.8670662807905248
    "AK167" = 1.212061784103416
    "AK168" = -1.186252107208071
    "AK169" = -0.1145011618341881
    "AK170" = 0.13960191438372
    "AK171" = 0.3583320470313792
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
